$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 708.6667
$ws.Range("I32").Value = 616.6667
$ws.Range("J32").Value = 800.6667
$ws.Range("K32").Value = 616.6667
$ws.Range("L32").Value = 800.6667
$ws.Range("M32").Value = -290.6667
$ws.Range("N32").Value = -1452.6667

$ws.Range("H46").Value = 4259.7144
$ws.Range("J46").Value = 4259.7144
$ws.Range("L46").Value = 12779.1432
$ws.Range("N46").Value = -13017.1432

$ws.Range("H60").Value = 4259.7144
$ws.Range("J60").Value = 4259.7144
$ws.Range("L60").Value = 12779.1432
$ws.Range("N60").Value = -13747.1432

$ws.Range("H76").Value = 4003144
$ws.Range("I76").Value = 4169816.8
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 4169816.8
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -4169501.8
$ws.Range("N76").Value = -3630

$ws.Range("H79").Value = 4003144
$ws.Range("I79").Value = 4169816.8
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 4169816.8
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -4168724.8
$ws.Range("N79").Value = -5184

$ws.Range("H86").Value = 66669864
$ws.Range("I86").Value = 4197.2
$ws.Range("J86").Value = 100002696
$ws.Range("K86").Value = 4197.2
$ws.Range("L86").Value = 100002696
$ws.Range("M86").Value = -3074.2
$ws.Range("N86").Value = -100004942

$ws.Range("H89").Value = 66669864
$ws.Range("I89").Value = 4197.2
$ws.Range("J89").Value = 100002696
$ws.Range("K89").Value = 20986
$ws.Range("L89").Value = 500013480
$ws.Range("M89").Value = -15370
$ws.Range("N89").Value = -500024712

$ws.Range("H132").Value = 1207.1296
$ws.Range("I132").Value = 1235.4318
$ws.Range("J132").Value = 1082.6
$ws.Range("K132").Value = 3706.2954
$ws.Range("L132").Value = 3247.8
$ws.Range("M132").Value = -1176.2954
$ws.Range("N132").Value = -8307.799999999999

$ws.Range("H137").Value = 773.4039
$ws.Range("I137").Value = 693.5806
$ws.Range("J137").Value = 891.2381
$ws.Range("K137").Value = 2080.7418
$ws.Range("L137").Value = 2673.7143
$ws.Range("M137").Value = 469.2582000000002
$ws.Range("N137").Value = -7773.7143

$ws.Range("H138").Value = 1609.4747
$ws.Range("I138").Value = 757.0492
$ws.Range("J138").Value = 2977.842
$ws.Range("K138").Value = 2271.1476
$ws.Range("L138").Value = 8933.526
$ws.Range("M138").Value = 2868.8524
$ws.Range("N138").Value = -19213.526

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1670.96
$ws.Range("I32").Value = 1621.3948
$ws.Range("J32").Value = 1827.9166
$ws.Range("K32").Value = 1621.3948
$ws.Range("L32").Value = 1827.9166
$ws.Range("M32").Value = -1334.3948
$ws.Range("N32").Value = -2401.9166

$ws.Range("H74").Value = 928.3570999999999
$ws.Range("I74").Value = 849.63635
$ws.Range("J74").Value = 1217
$ws.Range("K74").Value = 849.63635
$ws.Range("L74").Value = 1217
$ws.Range("M74").Value = 24.36365000000001
$ws.Range("N74").Value = -2965

$ws.Range("H77").Value = 928.3570999999999
$ws.Range("I77").Value = 849.63635
$ws.Range("J77").Value = 1217
$ws.Range("K77").Value = 4248.18175
$ws.Range("L77").Value = 6085
$ws.Range("M77").Value = 119.8182500000003
$ws.Range("N77").Value = -14821

$ws.Range("H102").Value = 2875.5
$ws.Range("I102").Value = 2344.75
$ws.Range("K102").Value = 2344.75
$ws.Range("M102").Value = -722.75

$ws.Range("H132").Value = 1025.4706
$ws.Range("I132").Value = 934.91113
$ws.Range("K132").Value = 2804.73339
$ws.Range("M132").Value = -274.7333899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4373.5557
$ws.Range("I105").Value = 3474.85
$ws.Range("J105").Value = 6941.2856
$ws.Range("K105").Value = 3474.85
$ws.Range("L105").Value = 6941.2856
$ws.Range("M105").Value = -1727.85
$ws.Range("N105").Value = -10435.2856

$ws.Range("H134").Value = 15963.072
$ws.Range("I134").Value = 1276.228
$ws.Range("J134").Value = 85725.586
$ws.Range("K134").Value = 3828.684
$ws.Range("L134").Value = 257176.758
$ws.Range("M134").Value = -1293.684
$ws.Range("N134").Value = -262246.758

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2488.377
$ws.Range("I31").Value = 2302.238
$ws.Range("J31").Value = 2899.842
$ws.Range("K31").Value = 2302.238
$ws.Range("L31").Value = 2899.842
$ws.Range("M31").Value = -2007.238
$ws.Range("N31").Value = -3489.842

$ws.Range("H34").Value = 2488.377
$ws.Range("I34").Value = 2302.238
$ws.Range("J34").Value = 2899.842
$ws.Range("K34").Value = 2302.238
$ws.Range("L34").Value = 2899.842
$ws.Range("M34").Value = -2100.238
$ws.Range("N34").Value = -3303.842

$ws.Range("H132").Value = 1309.8928
$ws.Range("I132").Value = 1011.1667
$ws.Range("J132").Value = 1847.6
$ws.Range("K132").Value = 3033.5001
$ws.Range("L132").Value = 5542.799999999999
$ws.Range("M132").Value = -503.5001000000002
$ws.Range("N132").Value = -10602.8

$ws.Range("H134").Value = 1173.4744
$ws.Range("I134").Value = 1095.8
$ws.Range("J134").Value = 1561.8462
$ws.Range("K134").Value = 3287.4
$ws.Range("L134").Value = 4685.5386
$ws.Range("M134").Value = -752.3999999999996
$ws.Range("N134").Value = -9755.5386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 920.2857
$ws.Range("I18").Value = 354.92307
$ws.Range("J18").Value = 1839
$ws.Range("K18").Value = 1064.76921
$ws.Range("L18").Value = 5517
$ws.Range("M18").Value = -895.7692099999999
$ws.Range("N18").Value = -5855

$ws.Range("H45").Value = 1031.6
$ws.Range("I45").Value = 460
$ws.Range("J45").Value = 1174.5
$ws.Range("K45").Value = 1380
$ws.Range("L45").Value = 3523.5
$ws.Range("M45").Value = -848
$ws.Range("N45").Value = -4587.5

$ws.Range("H131").Value = 6508060.5
$ws.Range("I131").Value = 167003260
$ws.Range("J131").Value = 1498.2162
$ws.Range("K131").Value = 501009780
$ws.Range("L131").Value = 4494.6486
$ws.Range("M131").Value = -501004740
$ws.Range("N131").Value = -14574.6486

$ws.Range("H137").Value = 31562.172
$ws.Range("I137").Value = 1825.9412
$ws.Range("J137").Value = 59646.39
$ws.Range("K137").Value = 5477.8236
$ws.Range("L137").Value = 178939.17
$ws.Range("M137").Value = -377.8235999999997
$ws.Range("N137").Value = -189139.17

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4350
$ws.Range("I70").Value = 4084.6155
$ws.Range("J70").Value = 4663.636
$ws.Range("K70").Value = 4084.6155
$ws.Range("L70").Value = 4663.636
$ws.Range("M70").Value = -3814.6155
$ws.Range("N70").Value = -5203.636

$ws.Range("H73").Value = 4350
$ws.Range("I73").Value = 4084.6155
$ws.Range("J73").Value = 4663.636
$ws.Range("K73").Value = 4084.6155
$ws.Range("L73").Value = 4663.636
$ws.Range("M73").Value = -3148.6155
$ws.Range("N73").Value = -6535.636

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1770.9153
$ws.Range("I132").Value = 1740.76
$ws.Range("J132").Value = 1938.4445
$ws.Range("K132").Value = 5222.28
$ws.Range("L132").Value = 5815.333500000001
$ws.Range("M132").Value = -2692.28
$ws.Range("N132").Value = -10875.3335

$ws.Range("H136").Value = 2526.1516
$ws.Range("I136").Value = 1348.6786
$ws.Range("J136").Value = 9120
$ws.Range("K136").Value = 4046.0358
$ws.Range("L136").Value = 27360
$ws.Range("M136").Value = -1496.0358
$ws.Range("N136").Value = -32460

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 736.0961
$ws.Range("I132").Value = 708.4186
$ws.Range("J132").Value = 868.3333
$ws.Range("K132").Value = 2125.2558
$ws.Range("L132").Value = 2604.9999
$ws.Range("M132").Value = 404.7442000000001
$ws.Range("N132").Value = -7664.9999

$ws.Range("H136").Value = 895.5789
$ws.Range("I136").Value = 889.6177
$ws.Range("J136").Value = 946.25
$ws.Range("K136").Value = 2668.8531
$ws.Range("L136").Value = 2838.75
$ws.Range("M136").Value = -118.8531000000003
$ws.Range("N136").Value = -7938.75
